$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these Price cells keep their literal text (avoid Excel auto-numeric conversion,
# which would silently drop significant trailing/represented digits, e.g. "5.00" -> 5).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"

# Apply the updated values scraped for this run.
$ws.Range("D2").Value = '43.161.32'
$ws.Range("E2").Value = '  -0.40%  '
$ws.Range("D3").Value = '2.354.32'
$ws.Range("E3").Value = '  +5.25%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '233.84'
$ws.Range("E5").Value = '  +1.56%  '
$ws.Range("D6").Value = '0.651'
$ws.Range("E6").Value = '  +1.91%  '
$ws.Range("D7").Value = '72.11'
$ws.Range("E7").Value = '  +14.14%  '
$ws.Range("E8").Value = '  -0.13%  '
$ws.Range("D9").Value = '0.497'
$ws.Range("E9").Value = '  +12.77%  '
$ws.Range("D10").Value = '0.0976'
$ws.Range("E10").Value = '  +2.30%  '
$ws.Range("D11").Value = '27.38'
$ws.Range("E11").Value = '  -0.56%  '
$ws.Range("E12").Value = '  +2.30%  '
$ws.Range("D13").Value = '2.704.10'
$ws.Range("E13").Value = '  +5.23%  '
$ws.Range("D14").Value = '16.20'
$ws.Range("E14").Value = '  +5.20%  '
$ws.Range("E15").Value = '  +4.58%  '
$ws.Range("D16").Value = '0.866'
$ws.Range("E16").Value = '  +5.09%  '
$ws.Range("D17").Value = '2.354.25'
$ws.Range("E17").Value = '  +4.87%  '
$ws.Range("D18").Value = '43.185.24'
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("D19").Value = '0.0₃0999'
$ws.Range("E19").Value = '  +3.87%  '
$ws.Range("D20").Value = '6.36'
$ws.Range("E20").Value = '  +4.81%  '
$ws.Range("D21").Value = '74.52'
$ws.Range("E21").Value = '  +2.46%  '
$ws.Range("D22").Value = '250.96'
$ws.Range("E22").Value = '  +2.11%  '
$ws.Range("D23").Value = '3.84'
$ws.Range("E23").Value = '  +4.95%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("D25").Value = '2.47'
$ws.Range("E25").Value = '  +2.18%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '10.06'
$ws.Range("E26").Value = '  +3.49%  '
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").Value = '2.21'
$ws.Range("E27").Value = '  -2.65%  '
$ws.Range("D28").Value = '22.45'
$ws.Range("E28").Value = '  +4.68%  '
$ws.Range("D29").Value = '172.35'
$ws.Range("E29").Value = '  -0.28%  '
$ws.Range("D30").Value = '1.54'
$ws.Range("E30").Value = '  +9.05%  '
$ws.Range("D31").Value = '0.133'
$ws.Range("E31").Value = '  +2.67%  '
$ws.Range("E32").Value = '  +2.65%  '
$ws.Range("D33").Value = '5.00'
$ws.Range("E33").Value = '  +2.83%  '
$ws.Range("D34").Value = '0.0692'
$ws.Range("E34").Value = '  +3.07%  '
$ws.Range("D35").Value = '5.07'
$ws.Range("E35").Value = '  +4.41%  '
$ws.Range("D36").Value = '3.74'
$ws.Range("E36").Value = '  +4.64%  '
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").Value = '2.43'
$ws.Range("E37").Value = '  +7.21%  '
$ws.Range("B38").Value = 'THORChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D38").Value = '6.52'
$ws.Range("E38").Value = '  +3.95%  '
$ws.Range("D39").Value = '0.0254'
$ws.Range("E39").Value = '  +2.01%  '
$ws.Range("D40").Value = '19.01'
$ws.Range("E40").Value = '  +12.05%  '
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("D42").Value = '8.90'
$ws.Range("E42").Value = '  +3.61%  '
$ws.Range("D43").Value = '4.51'
$ws.Range("E43").Value = '  +1.06%  '
$ws.Range("D44").Value = '99.20'
$ws.Range("E44").Value = '  +3.13%  '
$ws.Range("D45").Value = '1.16'
$ws.Range("E45").Value = '  +9.62%  '
$ws.Range("D46").Value = '0.0959'
$ws.Range("E46").Value = '  +2.17%  '
$ws.Range("E47").Value = '  +2.81%  '
$ws.Range("D48").Value = '1.442.25'
$ws.Range("E48").Value = '  -0.46%  '
$ws.Range("D49").Value = '2.578.36'
$ws.Range("E49").Value = '  +5.41%  '
$ws.Range("E50").Value = '  +0.88%  '
$ws.Range("E51").Value = '  -3.93%  '
